$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77; existing rows 77-186 shift down to 78-187.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new weekly price record.
$ws.Range("A77").Value = 9
$ws.Range("B77").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C77").Value = "Metropolitana"
$ws.Range("D77").Value = 44533
$ws.Range("E77").Value = 13
$ws.Range("F77").Value = 100112026
$ws.Range("G77").Value = "Haba"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 160
$ws.Range("K77").Value = 8000
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = 8500
$ws.Range("N77").Value = "$/saco 25 kilos"
$ws.Range("O77").Value = "Región del Maule"
$ws.Range("P77").Value = 340
$ws.Range("Q77").Value = 25
$ws.Range("R77").Value = "Hortaliza"

# Keep the date cell formatted the same way as the rest of column D.
$ws.Range("D77").NumberFormat = $ws.Range("D78").NumberFormat
